# Apply the score correction for Yuzvendra Chahal's two innings rows.
# Row 2: runs 1 -> 0, balls 3 -> 0
# Row 3: runs 0 -> 1, balls 0 -> 3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these as text (matching the original "numberStoredAsText" cells)
# rather than letting Excel auto-convert the numeric-looking strings to numbers.
$ws.Range("C2:D3").NumberFormat = "@"

$ws.Cells.Item(2, 3).Value = "0"
$ws.Cells.Item(2, 4).Value = "0"
$ws.Cells.Item(3, 3).Value = "1"
$ws.Cells.Item(3, 4).Value = "3"
